# Insert a new data row at row 100 (pushing the existing rows 100-223 down
# to 101-224), then populate the newly inserted row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("100:100").Insert()

$ws.Range("A100").Value = 10
$ws.Range("B100").Value = "Vega Modelo de Temuco"
$ws.Range("C100").Value = "La Araucanía"
$ws.Range("D100").Value = 44778
$ws.Range("E100").Value = 9
$ws.Range("F100").Value = 100112005
$ws.Range("G100").Value = "Puerro"
$ws.Range("H100").Value = "Azul de Maquehue"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 50
$ws.Range("K100").Value = 16000
$ws.Range("L100").Value = 16000
$ws.Range("M100").Value = 16000
$ws.Range("N100").Value = '$/docena de paquetes'
$ws.Range("O100").Value = "Provincia de Cautín"
$ws.Range("P100").Value = 1333
$ws.Range("Q100").Value = 12
$ws.Range("R100").Value = "Hortaliza"
